$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-5 per the commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 0
